$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's used range ends at row 472 (row 1 = header, data rows 2..472).
$lastRow = 472

# 1) Bump the "Förändrad" (last-changed) date in column C by two days
#    (45184 -> 45186) for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# 2) For the rows whose HYPERLINK() formulas in columns S, T, V, W, X, Y
#    only take a URL argument, add the case id (column A / "Beteckning")
#    as the second (friendly-name) argument.
$hyperlinkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    $nameCell = $ws.Cells.Item($r, 1)
    $name = $nameCell.Text

    if ([string]::IsNullOrEmpty($name)) {
        continue
    }

    foreach ($c in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $c)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula -match '^=HYPERLINK\("([^"]*)"\)$') {
            $url = $Matches[1]
            $cell.Formula = '=HYPERLINK("' + $url + '", "' + $name + '")'
        }
    }
}
